# Java Season 3 Notes and Programs till Sep 22nd
#
# Sheet1: tweak a couple of cells / add a couple of new ones.
# Sheet3:  replace the old scratch "emp id / ename" table with the FLM
#          automation config (browser / url / implicitWait / username /
#          password), bump the font size, and make it the active sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet1
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")

$ws1.Range("B2").Value = "reyaz"
$ws1.Range("D2").Value = "FLM"
$ws1.Range("A5").Value = "Passed"

# ---------------------------------------------------------------------
# Sheet3 - wipe the old content and lay down the new config table
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Cells.Clear() | Out-Null

$ws3.Range("A1").Value = "browser"
$ws3.Range("B1").Value = "chrome"
$ws3.Range("A2").Value = "url"
$ws3.Range("B2").Value = "https://adactinhotelapp.com/"
$ws3.Range("A3").Value = "implicitWait"
$ws3.Range("B3").Value = 30
$ws3.Range("A4").Value = "username"
$ws3.Range("B4").Value = "reyaz0806"
$ws3.Range("A5").Value = "password"
$ws3.Range("B5").Value = "reyaz123"

# Bigger font across the whole sheet (new style -> fontId 2, size 20)
$ws3.Cells.Font.Size = 20

# Nice readable column widths for the new key/value layout
$ws3.Columns.Item(1).AutoFit() | Out-Null
$ws3.Columns.Item(2).AutoFit() | Out-Null

# Make Sheet3 the active tab, with the same lingering selection the
# original author left behind.
$ws3.Activate() | Out-Null
$ws3.Range("B8").Select() | Out-Null
